$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Address (F2) was empty; fill it in — adds a new shared string "Anand,Vadodara"
$ws.Range("F2").Value = "Anand,Vadodara"

# IsBlackListed (J2) changes from TRUE to FALSE
$ws.Range("J2").Value = $false

# Move/update the active selection to J4 (was B4)
$ws.Range("J4").Select()
